{"js": "// Update the two-digit multiplication problems throughout the document.\n// Each old expression is unique within the document, so a direct\n// search-and-replace (preserving the run's formatting) is safe.\nconst replacements = [\n  [\"25\u00d746=\", \"79\u00d757=\"],\n  [\"69\u00d748=\", \"39\u00d747=\"],\n  [\"72\u00d723=\", \"97\u00d782=\"],\n  [\"29\u00d720=\", \"76\u00d719=\"],\n  [\"18\u00d747=\", \"17\u00d716=\"],\n  [\"19\u00d749=\", \"97\u00d723=\"],\n  [\"99\u00d726=\", \"82\u00d775=\"],\n  [\"42\u00d787=\", \"98\u00d742=\"],\n  [\"56\u00d774=\", \"21\u00d743=\"],\n  [\"26\u00d795=\", \"14\u00d738=\"],\n  [\"64\u00d798=\", \"12\u00d739=\"],\n  [\"62\u00d732=\", \"81\u00d760=\"],\n  [\"37\u00d754=\", \"39\u00d725=\"],\n  [\"76\u00d775=\", \"87\u00d784=\"],\n  [\"90\u00d727=\", \"19\u00d795=\"],\n  [\"37\u00d742=\", \"32\u00d722=\"],\n  [\"30\u00d738=\", \"67\u00d799=\"],\n  [\"61\u00d779=\", \"73\u00d753=\"],\n  [\"14\u00d752=\", \"21\u00d729=\"],\n  [\"57\u00d715=\", \"98\u00d790=\"],\n  [\"83\u00d735=\", \"52\u00d713=\"],\n  [\"88\u00d772=\", \"45\u00d790=\"],\n  [\"18\u00d778=\", \"79\u00d770=\"],\n  [\"42\u00d718=\", \"31\u00d723=\"],\n  [\"13\u00d736=\", \"26\u00d783=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit multiplication problems throughout the document.\n# Each old expression is unique within the document, so Find/Replace by\n# exact text (wildcards off) safely targets the correct run each time.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"25\u00d746=\", \"79\u00d757=\"),\n  @(\"69\u00d748=\", \"39\u00d747=\"),\n  @(\"72\u00d723=\", \"97\u00d782=\"),\n  @(\"29\u00d720=\", \"76\u00d719=\"),\n  @(\"18\u00d747=\", \"17\u00d716=\"),\n  @(\"19\u00d749=\", \"97\u00d723=\"),\n  @(\"99\u00d726=\", \"82\u00d775=\"),\n  @(\"42\u00d787=\", \"98\u00d742=\"),\n  @(\"56\u00d774=\", \"21\u00d743=\"),\n  @(\"26\u00d795=\", \"14\u00d738=\"),\n  @(\"64\u00d798=\", \"12\u00d739=\"),\n  @(\"62\u00d732=\", \"81\u00d760=\"),\n  @(\"37\u00d754=\", \"39\u00d725=\"),\n  @(\"76\u00d775=\", \"87\u00d784=\"),\n  @(\"90\u00d727=\", \"19\u00d795=\"),\n  @(\"37\u00d742=\", \"32\u00d722=\"),\n  @(\"30\u00d738=\", \"67\u00d799=\"),\n  @(\"61\u00d779=\", \"73\u00d753=\"),\n  @(\"14\u00d752=\", \"21\u00d729=\"),\n  @(\"57\u00d715=\", \"98\u00d790=\"),\n  @(\"83\u00d735=\", \"52\u00d713=\"),\n  @(\"88\u00d772=\", \"45\u00d790=\"),\n  @(\"18\u00d778=\", \"79\u00d770=\"),\n  @(\"42\u00d718=\", \"31\u00d723=\"),\n  @(\"13\u00d736=\", \"26\u00d783=\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
